# Update the DividendHistory sheet with the latest dividend entry.
# A new XD Date (19/03/2025) was declared, so insert a new row right
# below the header (row 2) and push the existing history down, then
# populate the new row with the XD Date, Pay Date and Gross Dividend.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DividendHistory")

# Insert a new blank row at row 2 (shifts existing rows 2.. down by one)
$ws.Rows("2:2").Insert()

# Populate the newly inserted row with the new dividend record.
$ws.Range("A2").Value = "19/03/2025"
$ws.Range("B2").Value = "19/03/2025"
$ws.Range("C2").Value = "0.005"
